$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells whose new value looks like a plain decimal number
# (single dot, e.g. "111.36") stay stored as text, matching the rest of
# column D (prices with thousands separators, like "49.330.64", are never
# auto-converted since Excel cannot parse them as numbers).
$textRows = @(5,6,7,10,11,14,17,19,23,24,26,28,32,36,37,39,40,41,47,49,50,51)
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "49.330.64"
$ws.Range("E2").Value = "  -1.08%  "
$ws.Range("D3").Value = "2.625.03"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "111.36"
$ws.Range("E5").Value = "  -2.44%  "
$ws.Range("D6").Value = "325.41"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "0.524"
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").Value = "39.35"
$ws.Range("E10").Value = "  -4.06%  "
$ws.Range("D11").Value = "19.99"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("E12").Value = "  -1.29%  "
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").Value = "7.54"
$ws.Range("E14").Value = "  +2.54%  "
$ws.Range("D15").Value = "3.041.36"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "2.632.85"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "0.850"
$ws.Range("E17").Value = "  -2.13%  "
$ws.Range("D18").Value = "49.355.18"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "13.20"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("D23").Value = "267.63"
$ws.Range("E23").Value = "  -3.49%  "
$ws.Range("D24").Value = "68.91"
$ws.Range("E24").Value = "  -4.36%  "
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").Value = "25.97"
$ws.Range("E26").Value = "  -3.01%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "10.14"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("E31").Value = "  -4.56%  "
$ws.Range("D32").Value = "49.46"
$ws.Range("E32").Value = "  -1.69%  "
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "18.97"
$ws.Range("E36").Value = "  -2.95%  "
$ws.Range("D37").Value = "4.93"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("E38").Value = "  -2.43%  "
$ws.Range("D39").Value = "3.08"
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").Value = "128.80"
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("D41").Value = "22.64"
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("E43").Value = "  -0.88%  "
$ws.Range("E44").Value = "  +5.99%  "
$ws.Range("D45").Value = "2.057.48"
$ws.Range("E45").Value = "  -1.08%  "
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("D47").Value = "2.12"
$ws.Range("E47").Value = "  +7.64%  "
$ws.Range("E48").Value = "  -5.98%  "
$ws.Range("D49").Value = "8.84"
$ws.Range("E49").Value = "  -2.93%  "
$ws.Range("D50").Value = "5.19"
$ws.Range("E50").Value = "  -3.41%  "
$ws.Range("D51").Value = "58.33"
$ws.Range("E51").Value = "  -3.23%  "
